$wb = $excel.ActiveWorkbook

# --- Update the shared "Ready for handoff" status text to "Handback transform
# failed" everywhere it is used (it's a single shared string reused by the
# Overview summary row and the per-locale "Status" column for the same
# c3f54ab4-... record), so the underlying shared string is edited in place
# rather than duplicated.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handback transform failed"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: widen the "Error Detail" column (P) and populate the
# handback-mismatch error message for row 3 (c3f54ab4-... record).
$wsZh.Columns.Item(16).ColumnWidth = 39.17
$wsZh.Range("P3").Value = "Handback file name: wb1y5eah.b4l is different with handoff file name: c3f54ab4-bc39-48dd-9ef4-adc709e15aac.056be68b88399324172ebb9f70b4ada65403fe0c.zh-cn."

# --- de-de sheet: widen the "Error Detail" column (P) and populate the
# handback-mismatch error message for row 3 (c3f54ab4-... record).
$wsDe.Columns.Item(16).ColumnWidth = 39.17
$wsDe.Range("P3").Value = "Handback file name: wb1y5eah.b4l is different with handoff file name: c3f54ab4-bc39-48dd-9ef4-adc709e15aac.056be68b88399324172ebb9f70b4ada65403fe0c.de-de."
